$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# New rows of data to append starting at row 42 (columns: A, C, F).
# All values are text (matching the original sheet's t="str" convention),
# so NumberFormat is forced to "@" (Text) before assignment to avoid
# Excel auto-converting numeric-looking strings into numbers.
$newRows = @(
    @{ Row = 42; A = $null; C = "106_绣球单瓣粉_Hydrangea Pink S_Hydrangea L._1stem"; F = "25" },
    @{ Row = 43; A = $null; C = "345_天竺少女_Cryptomeria`nKashiwaba_undefined_1bunch"; F = "15" },
    @{ Row = 44; A = $null; C = "328_卢荀草_undefined_undefined_1bunch"; F = "25" },
    @{ Row = 45; A = $null; C = "540_糖棉_gomphocarpus fruticosus_undefined_1bunch"; F = "6" },
    @{ Row = 46; A = $null; C = "462_五针松_undefined_undefined_1bunch"; F = "5" },
    @{ Row = 47; A = $null; C = "341_南天竹绿_undefined_Nandina domestica Thunb._1bunch"; F = "5" },
    @{ Row = 48; A = $null; C = "348_万年青_undefined_undefined_1bunch"; F = "15" },
    @{ Row = 49; A = "8";   C = "316_尤加利叶大叶_Eucalyptus Cinerea_undefined_1bunch"; F = "12" },
    @{ Row = 50; A = $null; C = "317_尤加利叶细叶_Eucalyptus Parvifolia_undefined_1bunch"; F = "10" },
    @{ Row = 51; A = $null; C = "319_尤加利叶带果_Eucalyptus leaves with small pods_undefined_1bunch"; F = $null }
)

foreach ($row in $newRows) {
    $r = $row.Row
    if ($row.A -ne $null) {
        $cell = $ws.Cells.Item($r, 1)
        $cell.NumberFormat = "@"
        $cell.Value = $row.A
    }
    if ($row.C -ne $null) {
        $cell = $ws.Cells.Item($r, 3)
        $cell.NumberFormat = "@"
        $cell.Value = $row.C
    }
    if ($row.F -ne $null) {
        $cell = $ws.Cells.Item($r, 6)
        $cell.NumberFormat = "@"
        $cell.Value = $row.F
    }
}

# Update Summary sheet G2 value (concatenated Number column string got
# extended with the new rows' values).
$wsSummary = $wb.Worksheets.Item("Summary")
$cellG2 = $wsSummary.Range("G2")
$cellG2.NumberFormat = "@"
$cellG2.Value = "010136731028153831510192518411810104101978125302525251540605061013122515256551512100"
